$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserPasswords")

# Update the password hash stored in B6 to match the common hash value
# used across the rest of the sheet.
$ws.Range("B6").Value = "5e884898da28047151d0e56f8dc6292773603d0d6aabbdd62a11ef721d1542d8"
